$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '29.404.10'
$ws.Cells.Item(2, 5).Value = '  -0.26%  '
$ws.Cells.Item(3, 4).Value = '1.841.65'
$ws.Cells.Item(3, 5).Value = '  -0.38%  '
$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = '0.9996'
$ws.Cells.Item(4, 5).Value = '  +0.14%  '
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '239.24'
$ws.Cells.Item(5, 5).Value = '  -0.45%  '
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '0.6278'
$ws.Cells.Item(6, 5).Value = '  -0.26%  '
$ws.Cells.Item(7, 5).Value = '  +0.16%  '
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = '0.07399'
$ws.Cells.Item(8, 5).Value = '  -1.16%  '
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '0.2892'
$ws.Cells.Item(9, 5).Value = '  -0.61%  '
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '24.89'
$ws.Cells.Item(10, 5).Value = '  +1.12%  '
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '0.07723'
$ws.Cells.Item(11, 5).Value = '  -0.26%  '
$ws.Cells.Item(12, 4).Value = '1.847.92'
$ws.Cells.Item(12, 5).Value = '  -0.30%  '
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '4.963'
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '0.6692'
$ws.Cells.Item(14, 5).Value = '  -1.74%  '
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '0.00001037'
$ws.Cells.Item(15, 5).Value = '  -0.60%  '
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '81.66'
$ws.Cells.Item(16, 5).Value = '  -0.63%  '
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = '6.247'
$ws.Cells.Item(17, 5).Value = '  +0.04%  '
$ws.Cells.Item(18, 4).Value = '29.433.53'
$ws.Cells.Item(18, 5).Value = '  -0.16%  '
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '233.99'
$ws.Cells.Item(19, 5).Value = '  +1.95%  '
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '12.29'
$ws.Cells.Item(20, 5).Value = '  -0.84%  '
$ws.Cells.Item(21, 5).Value = '  +0.13%  '
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '7.295'
$ws.Cells.Item(22, 5).Value = '  -3.29%  '
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '1.001'
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '156.97'
$ws.Cells.Item(24, 5).Value = '  -1.56%  '
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '8.472'
$ws.Cells.Item(25, 5).Value = '  -0.66%  '
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '0.1343'
$ws.Cells.Item(26, 5).Value = '  -1.91%  '
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '17.31'
$ws.Cells.Item(27, 5).Value = '  -1.31%  '
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '0.07225'
$ws.Cells.Item(28, 5).Value = '  +10.39%  '
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '1.487'
$ws.Cells.Item(29, 5).Value = '  +4.77%  '
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '1.482'
$ws.Cells.Item(30, 5).Value = '  -0.33%  '
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '4.028'
$ws.Cells.Item(31, 5).Value = '  -1.81%  '
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '4.030'
$ws.Cells.Item(32, 5).Value = '  -1.89%  '
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '1.159'
$ws.Cells.Item(33, 5).Value = '  +1.12%  '
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '1.817'
$ws.Cells.Item(34, 5).Value = '  -1.10%  '
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '0.7089'
$ws.Cells.Item(35, 5).Value = '  +1.35%  '
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '2.582'
$ws.Cells.Item(36, 5).Value = '  +0.18%  '
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '0.01835'
$ws.Cells.Item(37, 5).Value = '  -1.50%  '
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '2.794'
$ws.Cells.Item(38, 5).Value = '  -1.71%  '
$ws.Cells.Item(39, 4).Value = '1.232.47'
$ws.Cells.Item(39, 5).Value = '  -2.66%  '
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '6.784'
$ws.Cells.Item(40, 5).Value = '  -0.68%  '
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '0.9519'
$ws.Cells.Item(41, 5).Value = '  +1.87%  '
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '1.001'
$ws.Cells.Item(42, 5).Value = '  +0.10%  '
$ws.Cells.Item(43, 4).Value = '2.005.46'
$ws.Cells.Item(43, 5).Value = '  -0.51%  '
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '101.15'
$ws.Cells.Item(44, 5).Value = '  -0.20%  '
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '65.29'
$ws.Cells.Item(45, 5).Value = '  -1.51%  '
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '0.00000000117'
$ws.Cells.Item(46, 5).Value = '  -1.36%  '
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '1.697'
$ws.Cells.Item(47, 5).Value = '  -2.18%  '
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '6.947'
$ws.Cells.Item(48, 5).Value = '  -2.19%  '
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '8.903'
$ws.Cells.Item(49, 5).Value = '  -1.13%  '
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '0.1131'
$ws.Cells.Item(50, 5).Value = '  -2.92%  '
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '0.3884'
$ws.Cells.Item(51, 5).Value = '  -1.93%  '
